$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = 43866

$ws.Range("B14").Value = "Work with numba on Windows 10"
$ws.Range("B15").Value = "Try to install pyopencl without import error"
$ws.Range("B16").Value = "Afterwards, attempt to use llspy and spimagine, which both have opencl dependencies"
$ws.Range("B17").Value = "Cannot install pyopencl without import error - give up for now and look to working with Linux"

$ws.Range("B17").Select()
